$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "is_locked_lbl" / "is_enabled_lbl" template columns (the
# inlineForeignTabs editable aggregate table no longer needs per-row
# lock/enable dropdown validation columns). The remaining "order_by" and
# "rem" headers shift left into C1/D1, and the old E1/F1 cells go away.
$ws.Range("C1").Value = "<%=comment.order_by%>"
$ws.Range("D1").Value = "<%=comment.rem%>"
$ws.Range("E1:F1").ClearContents()
